$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 69
$ws.Range("H69").Value = 5989
$ws.Range("I69").Value = 5989
$ws.Range("K69").Value = 17967
$ws.Range("M69").Value = -17093
# Row 72
$ws.Range("H72").Value = 5989
$ws.Range("I72").Value = 5989
$ws.Range("K72").Value = 53901
$ws.Range("M72").Value = -49533
# Row 107
$ws.Range("H107").Value = 33930116
$ws.Range("I107").Value = 15626206
$ws.Range("J107").Value = 58335332
$ws.Range("K107").Value = 15626206
$ws.Range("L107").Value = 58335332
$ws.Range("M107").Value = -15624286
$ws.Range("N107").Value = -58339172
# Row 113
$ws.Range("H113").Value = 12355104
$ws.Range("I113").Value = 37038908
$ws.Range("J113").Value = 13202
$ws.Range("K113").Value = 37038908
$ws.Range("L113").Value = 13202
$ws.Range("M113").Value = -37035654
$ws.Range("N113").Value = -19710
# Row 132
$ws.Range("H132").Value = 1940.2
$ws.Range("I132").Value = 1548.48
$ws.Range("K132").Value = 4645.440000000001
$ws.Range("M132").Value = -2115.440000000001
# Row 137
$ws.Range("H137").Value = 3170.4285
$ws.Range("J137").Value = 2849.1667
$ws.Range("L137").Value = 8547.500100000001
$ws.Range("N137").Value = -13647.5001

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2079.0605
$ws.Range("I2").Value = 1312.1177
$ws.Range("J2").Value = 2893.9375
$ws.Range("K2").Value = 1312.1177
$ws.Range("L2").Value = 2893.9375
$ws.Range("M2").Value = -1199.1177
$ws.Range("N2").Value = -3119.9375
# Row 32
$ws.Range("H32").Value = 2235054.2
$ws.Range("I32").Value = 2317463.5
$ws.Range("K32").Value = 2317463.5
$ws.Range("M32").Value = -2317176.5
# Row 41
$ws.Range("H41").Value = 585.3333
$ws.Range("I41").Value = 585.3333
$ws.Range("K41").Value = 585.3333
$ws.Range("M41").Value = -171.3333
# Row 45
$ws.Range("H45").Value = 6460.3125
$ws.Range("I45").Value = 2548
$ws.Range("K45").Value = 2548
$ws.Range("M45").Value = -2171
# Row 110
$ws.Range("I110").Value = 1467.5714
$ws.Range("J110").Value = 333333340
$ws.Range("K110").Value = 1467.5714
$ws.Range("L110").Value = 333333340
$ws.Range("M110").Value = 577.4286
$ws.Range("N110").Value = -333337430
# Row 116
$ws.Range("H116").Value = 2079.0605
$ws.Range("I116").Value = 1312.1177
$ws.Range("J116").Value = 2893.9375
$ws.Range("K116").Value = 1312.1177
$ws.Range("L116").Value = 2893.9375
$ws.Range("M116").Value = 981.8823
$ws.Range("N116").Value = -7481.9375
# Row 122
$ws.Range("H122").Value = 30060.875
$ws.Range("I122").Value = 51122
$ws.Range("K122").Value = 153366
$ws.Range("M122").Value = -150916

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2079.0605
$ws.Range("I3").Value = 1312.1177
$ws.Range("J3").Value = 2893.9375
$ws.Range("K3").Value = 1312.1177
$ws.Range("L3").Value = 2893.9375
$ws.Range("M3").Value = -1198.1177
$ws.Range("N3").Value = -3121.9375
# Row 54
$ws.Range("H54").Value = 1083
$ws.Range("I54").Value = 1083
$ws.Range("K54").Value = 1083
$ws.Range("M54").Value = -599
# Row 96
$ws.Range("H96").Value = 27097.857
$ws.Range("I96").Value = 16770
$ws.Range("J96").Value = 52917.5
$ws.Range("K96").Value = 16770
$ws.Range("L96").Value = 52917.5
$ws.Range("M96").Value = -14024
$ws.Range("N96").Value = -58409.5
# Row 99
$ws.Range("H99").Value = 3498936.2
$ws.Range("I99").Value = 2098.9285
$ws.Range("K99").Value = 2098.9285
$ws.Range("M99").Value = -600.9285
# Row 102
$ws.Range("H102").Value = 21313.334
$ws.Range("I102").Value = 4999
$ws.Range("J102").Value = 53942
$ws.Range("K102").Value = 4999
$ws.Range("L102").Value = 53942
$ws.Range("M102").Value = -1754
$ws.Range("N102").Value = -60432
# Row 105
$ws.Range("H105").Value = 3861.625
$ws.Range("I105").Value = 2999.75
$ws.Range("K105").Value = 2999.75
$ws.Range("M105").Value = -1252.75
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9330.058000000001
$ws.Range("J31").Value = 14361.05
$ws.Range("L31").Value = 14361.05
$ws.Range("N31").Value = -14951.05
# Row 34
$ws.Range("H34").Value = 9330.058000000001
$ws.Range("J34").Value = 14361.05
$ws.Range("L34").Value = 14361.05
$ws.Range("N34").Value = -14765.05
# Row 98
$ws.Range("H98").Value = 25000
$ws.Range("I98").Value = 25000
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 25000
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -22754
$ws.Range("N98").ClearContents()
# Row 122
$ws.Range("H122").Value = 74061.36
$ws.Range("I122").Value = 1216
$ws.Range("J122").Value = 146906.72
$ws.Range("K122").Value = 3648
$ws.Range("L122").Value = 440720.16
$ws.Range("M122").Value = -1198
$ws.Range("N122").Value = -445620.16
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 656.6
$ws.Range("I12").Value = 193.75
$ws.Range("J12").Value = 824.9091
$ws.Range("K12").Value = 581.25
$ws.Range("L12").Value = 2474.7273
$ws.Range("M12").Value = -408.25
$ws.Range("N12").Value = -2820.7273
# Row 98
$ws.Range("H98").Value = 943.25
$ws.Range("I98").Value = 736.6667
$ws.Range("J98").Value = 1149.8334
$ws.Range("K98").Value = 2210.0001
$ws.Range("L98").Value = 3449.5002
$ws.Range("M98").Value = -712.0001000000002
$ws.Range("N98").Value = -6445.5002
# Row 132
$ws.Range("H132").Value = 9383.963
$ws.Range("I132").Value = 3885.5715
$ws.Range("K132").Value = 34970.1435
$ws.Range("M132").Value = -32440.1435

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 158.54546
$ws.Range("J2").Value = 99.5
$ws.Range("L2").Value = 99.5
$ws.Range("N2").Value = -325.5
# Row 11
$ws.Range("H11").Value = 77363640
$ws.Range("I11").Value = 5100002
$ws.Range("K11").Value = 5100002
$ws.Range("M11").Value = -5099863
# Row 102
$ws.Range("H102").Value = 4975.815
$ws.Range("I102").Value = 4015.3635
$ws.Range("K102").Value = 4015.3635
$ws.Range("M102").Value = -2393.3635
# Row 121
$ws.Range("H121").Value = 49663.668
$ws.Range("J121").Value = 49663.668
$ws.Range("L121").Value = 49663.668
$ws.Range("N121").Value = -53157.668
# Row 122
$ws.Range("H122").Value = 3153393.5
$ws.Range("I122").Value = 6039213.5
$ws.Range("J122").Value = 5226.091
$ws.Range("K122").Value = 18117640.5
$ws.Range("L122").Value = 15678.273
$ws.Range("M122").Value = -18115190.5
$ws.Range("N122").Value = -20578.273
# Row 132
$ws.Range("H132").Value = 5595.1113
$ws.Range("I132").Value = 2183.7368
$ws.Range("K132").Value = 6551.2104
$ws.Range("M132").Value = -4021.2104

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 6968.6313
$ws.Range("J40").Value = 7877.231
$ws.Range("L40").Value = 7877.231
$ws.Range("N40").Value = -8149.231

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 56
$ws.Range("H56").Value = 34957
$ws.Range("J56").Value = 34957
$ws.Range("L56").Value = 34957
$ws.Range("N56").Value = -36385
# Row 100
$ws.Range("H100").Value = 874.6875
$ws.Range("I100").Value = 545.4545000000001
$ws.Range("K100").Value = 1090.909
$ws.Range("M100").Value = -549.9090000000001
# Row 107
$ws.Range("H107").Value = 886.5263
$ws.Range("I107").Value = 761.4286
$ws.Range("K107").Value = 2284.2858
$ws.Range("M107").Value = -364.2857999999997
# Row 122
$ws.Range("H122").Value = 121342.2
$ws.Range("I122").Value = 218861.8
$ws.Range("K122").Value = 656585.3999999999
$ws.Range("M122").Value = -654135.3999999999
# Row 136
$ws.Range("H136").Value = 32295444
$ws.Range("I136").Value = 100001290
$ws.Range("J136").Value = 54565.24
$ws.Range("K136").Value = 300003870
$ws.Range("L136").Value = 163695.72
$ws.Range("M136").Value = -300001320
$ws.Range("N136").Value = -168795.72
